# Site rebuild: drop the trailing "Ver no Jupiter / Salvar em pdf / Salvar
# em docx" + "© 2020 ... Jekyll ..." footer block (plus the blank paragraph
# right before it) that used to follow the final "LOQ4053: Balanços de
# Massa e Energia (Requisito fraco)" requirement line.

$d = $word.ActiveDocument

$startMarker = "LOQ4053: Balanços de Massa e Energia (Requisito fraco)"
$endMarker = "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

# Find the paragraph that ends the LOQ4053 requirement line -- deletion
# should start right after it, so its own paragraph mark is kept.
$anchorPara = $null
# Find the paragraph holding the copyright/footer line -- deletion should
# run through the end of this paragraph, removing its paragraph mark too.
$footerPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($anchorPara -eq $null -and $t -match [regex]::Escape($startMarker)) {
        $anchorPara = $p
    }
    if ($footerPara -eq $null -and $t -match [regex]::Escape($endMarker)) {
        $footerPara = $p
    }
}

if ($anchorPara -ne $null -and $footerPara -ne $null) {
    # Removes: the blank paragraph + the "Ver no Jupiter..." paragraph +
    # the "© 2020..." paragraph, leaving the LOQ4053 paragraph directly
    # followed by whatever paragraph used to come after the footer block.
    $toDelete = $d.Range($anchorPara.Range.End, $footerPara.Range.End)
    $toDelete.Delete()
}
